# Add a new "statement_mappings" worksheet at the end of the workbook,
# populated with the statement -> statement_type mapping table, and the
# same bold/bordered/centered header style used by the other mapping
# sheets in this workbook.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Add()
$ws.Name = "statement_mappings"

# Header row
$ws.Range("A1").Value = "statement"
$ws.Range("B1").Value = "statement_type"
$ws.Range("A1:B1").Font.Bold = $true
$ws.Range("A1:B1").HorizontalAlignment = -4108
$ws.Range("A1:B1").VerticalAlignment = -4160
$ws.Range("A1:B1").Borders.LineStyle = 1

# Data rows: Oracle PL/SQL statement keyword -> mapped statement type
$data = @(
    @("SELECT", "select_statement"),
    @("INSERT", "insert_statement"),
    @("UPDATE", "update_statement"),
    @("DELETE", "delete_statement"),
    @("RAISE", "raise_statement"),
    @("NULL ", "null_statement"),
    @("RETURN", "return_statement")
)

for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 1).Value = $data[$i][0]
    $ws.Cells.Item($row, 2).Value = $data[$i][1]
}

# Move the newly added sheet (Excel inserts new sheets at the front) to
# become the last tab, after "function_mappings".
$ws.Move($null, $wb.Worksheets.Item($wb.Worksheets.Count))
